# Implement NORMALITY output file refactor:
#  - "Vowels" sheet: unify the number format of the D4:K6 data block to a
#    single plain one-decimal format ("0.0") instead of the old mix of
#    percentage ("0.0%") on rows 4-5 and plain decimal on row 6.
#  - "Cons manner" sheet: the stray wordlist-length numbers that had leaked
#    into column E (rows 4-6) are removed; that data doesn't belong on this
#    sheet.

$wb = $excel.ActiveWorkbook

$vowels = $wb.Worksheets.Item("Vowels")
$vowels.Range("D4:K6").NumberFormat = "0.0"

$consManner = $wb.Worksheets.Item("Cons manner")
$consManner.Range("E4:E6").ClearContents()
